$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 142-143; this shifts the former rows 142-222 down to 144-224
$ws.Rows("142:143").Insert()

# New row 142 data
$ws.Cells.Item(142, 1).Value2  = 9
$ws.Cells.Item(142, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(142, 3).Value2  = "Metropolitana"
$ws.Cells.Item(142, 4).Value2  = 44680
$ws.Cells.Item(142, 5).Value2  = 13
$ws.Cells.Item(142, 6).Value2  = 100112003
$ws.Cells.Item(142, 7).Value2  = "Ajo"
$ws.Cells.Item(142, 8).Value2  = "Chino"
$ws.Cells.Item(142, 9).Value2  = "1a (guarda)"
$ws.Cells.Item(142, 10).Value2 = 106
$ws.Cells.Item(142, 11).Value2 = 7000
$ws.Cells.Item(142, 12).Value2 = 7000
$ws.Cells.Item(142, 13).Value2 = 7000
$ws.Cells.Item(142, 14).Value2 = "`$/trenza 50 unidades"
$ws.Cells.Item(142, 15).Value2 = "Provincia de Talagante"
$ws.Cells.Item(142, 16).Value2 = 1400
$ws.Cells.Item(142, 17).Value2 = 5
$ws.Cells.Item(142, 18).Value2 = "Hortaliza"

# New row 143 data
$ws.Cells.Item(143, 1).Value2  = 9
$ws.Cells.Item(143, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(143, 3).Value2  = "Metropolitana"
$ws.Cells.Item(143, 4).Value2  = 44680
$ws.Cells.Item(143, 5).Value2  = 13
$ws.Cells.Item(143, 6).Value2  = 100112003
$ws.Cells.Item(143, 7).Value2  = "Ajo"
$ws.Cells.Item(143, 8).Value2  = "Chino"
$ws.Cells.Item(143, 9).Value2  = "2a (guarda)"
$ws.Cells.Item(143, 10).Value2 = 61
$ws.Cells.Item(143, 11).Value2 = 6000
$ws.Cells.Item(143, 12).Value2 = 6000
$ws.Cells.Item(143, 13).Value2 = 6000
$ws.Cells.Item(143, 14).Value2 = "`$/trenza 50 unidades"
$ws.Cells.Item(143, 15).Value2 = "Provincia de Talagante"
$ws.Cells.Item(143, 16).Value2 = 1200
$ws.Cells.Item(143, 17).Value2 = 5
$ws.Cells.Item(143, 18).Value2 = "Hortaliza"

# Ensure the date column keeps the same number format as the rest of column D
$ws.Range("D142:D143").NumberFormat = $ws.Range("D144").NumberFormat
